# Updates cryptos list values (price + 1h volume change) per the
# upstream GitHub Actions scrape commit. A handful of "Price" cells hold
# purely-numeric-looking strings (e.g. "1.00", "9.59") that must stay as
# literal text (matching the source inlineStr cells), so those are forced
# to Text format before assignment; everything else (values containing a
# second "." like "41.108.12", the unicode subscript in "0.0₃0926", the
# "  +/-x.xx%  " volume strings, coin names and links) is never
# auto-numeric and can be assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- cells that need forced text format (number-like strings) ----
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "316.96"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "89.65"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.498"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0834"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "31.98"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.88"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.773"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "71.50"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "235.12"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.88"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "24.17"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.59"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "34.86"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "155.90"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0749"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.51"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.97"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "16.67"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.92"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.26"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "18.84"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0276"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.90"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.54"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "95.18"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "73.88"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "52.19"

# ---- plain text assignments ----
$ws.Cells.Item(2, 4).Value = "41.108.12"
$ws.Cells.Item(2, 5).Value = "  -1.23%  "
$ws.Cells.Item(3, 4).Value = "2.427.19"
$ws.Cells.Item(3, 5).Value = "  -1.85%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 5).Value = "  -0.27%  "
$ws.Cells.Item(6, 5).Value = "  -2.79%  "
$ws.Cells.Item(7, 5).Value = "  -2.60%  "
$ws.Cells.Item(8, 5).Value = "  +0.07%  "
$ws.Cells.Item(9, 5).Value = "  -3.35%  "
$ws.Cells.Item(10, 5).Value = "  -2.76%  "
$ws.Cells.Item(11, 5).Value = "  -3.37%  "
$ws.Cells.Item(12, 5).Value = "  -2.16%  "
$ws.Cells.Item(13, 4).Value = "2.798.08"
$ws.Cells.Item(13, 5).Value = "  -1.99%  "
$ws.Cells.Item(14, 5).Value = "  -2.38%  "
$ws.Cells.Item(15, 5).Value = "  +1.57%  "
$ws.Cells.Item(16, 4).Value = "2.424.29"
$ws.Cells.Item(16, 5).Value = "  -2.24%  "
$ws.Cells.Item(17, 5).Value = "  -2.16%  "
$ws.Cells.Item(18, 4).Value = "41.030.24"
$ws.Cells.Item(18, 5).Value = "  -1.31%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0926"
$ws.Cells.Item(19, 5).Value = "  -2.69%  "
$ws.Cells.Item(20, 5).Value = "  -3.57%  "
$ws.Cells.Item(21, 5).Value = "  +0.28%  "
$ws.Cells.Item(22, 5).Value = "  -2.61%  "
$ws.Cells.Item(23, 5).Value = "  -2.00%  "
$ws.Cells.Item(24, 5).Value = "  -2.20%  "
$ws.Cells.Item(25, 5).Value = "  +0.06%  "
$ws.Cells.Item(26, 5).Value = "  -2.18%  "
$ws.Cells.Item(27, 5).Value = "  -2.21%  "
$ws.Cells.Item(28, 5).Value = "  -3.01%  "
$ws.Cells.Item(29, 5).Value = "  -2.81%  "
$ws.Cells.Item(30, 5).Value = "  -3.85%  "
$ws.Cells.Item(31, 5).Value = "  -2.78%  "
$ws.Cells.Item(32, 5).Value = "  -4.50%  "
$ws.Cells.Item(33, 5).Value = "  +0.02%  "
$ws.Cells.Item(34, 2).Value = "Hedera"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(34, 5).Value = "  -2.86%  "
$ws.Cells.Item(35, 2).Value = "WEMIXToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(35, 5).Value = "  -2.59%  "
$ws.Cells.Item(36, 5).Value = "  +2.17%  "
$ws.Cells.Item(37, 5).Value = "  -3.35%  "
$ws.Cells.Item(38, 5).Value = "  -1.06%  "
$ws.Cells.Item(39, 5).Value = "  -2.22%  "
$ws.Cells.Item(40, 5).Value = "  -2.23%  "
$ws.Cells.Item(41, 5).Value = "  -1.25%  "
$ws.Cells.Item(42, 4).Value = "1.997.45"
$ws.Cells.Item(43, 5).Value = "  -8.67%  "
$ws.Cells.Item(44, 5).Value = "  -1.95%  "
$ws.Cells.Item(45, 5).Value = "  -3.56%  "
$ws.Cells.Item(46, 5).Value = "  -2.84%  "
$ws.Cells.Item(47, 5).Value = "  +3.66%  "
$ws.Cells.Item(48, 4).Value = "2.654.35"
$ws.Cells.Item(48, 5).Value = "  -2.13%  "
$ws.Cells.Item(49, 5).Value = "  -2.46%  "
$ws.Cells.Item(50, 5).Value = "  +0.42%  "
$ws.Cells.Item(51, 5).Value = "  -0.38%  "
